$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in A1 (st12345 -> st13245)
$ws.Range("A1").Value = "st13245"

# New header cell B1 - copy formatting (style) from A1, then set its value
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "st1"

# New data cells in column B (plain, unstyled, text)
$ws.Range("B2").Value = "d"
$ws.Range("B3").Value = "d"

# A4 corrected from 5 to 0, keep it stored as text (not a number)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "0"
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

# New average/listing output in B4 (plain, unstyled, text)
$ws.Range("B4").Value = "3, 5, 5, 5"
